$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 90.210526
$ws.Range("I11").Value = 90.210526
$ws.Range("K11").Value = 90.210526
$ws.Range("M11").Value = 49.789474

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3029.9
$ws.Range("I94").Value = 3029.9
$ws.Range("K94").Value = 3029.9
$ws.Range("M94").Value = -2578.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 10527
$ws.Range("I98").Value = 11125.9375
$ws.Range("K98").Value = 11125.9375
$ws.Range("M98").Value = -9627.9375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 66667028
$ws.Range("J111").Value = 111111496
$ws.Range("L111").Value = 333334488
$ws.Range("N111").Value = -333340622

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4498.5
$ws.Range("J116").Value = 3997
$ws.Range("L116").Value = 3997
$ws.Range("N116").Value = -10881

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 10527
$ws.Range("I122").Value = 11125.9375
$ws.Range("K122").Value = 33377.8125
$ws.Range("M122").Value = -30927.8125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1799.6364
$ws.Range("I129").Value = 799.5
$ws.Range("K129").Value = 2398.5
$ws.Range("M129").Value = 2601.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2184.9375
$ws.Range("I2").Value = 895.9
$ws.Range("J2").Value = 4333.3335
$ws.Range("K2").Value = 895.9
$ws.Range("L2").Value = 4333.3335
$ws.Range("M2").Value = -782.9
$ws.Range("N2").Value = -4559.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3245.5122
$ws.Range("I32").Value = 3159.6843
$ws.Range("K32").Value = 3159.6843
$ws.Range("M32").Value = -2872.6843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1767
$ws.Range("I45").Value = 1644.4
$ws.Range("K45").Value = 1644.4
$ws.Range("M45").Value = -1267.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4395.7666
$ws.Range("I61").Value = 4034.8462
$ws.Range("K61").Value = 4034.8462
$ws.Range("M61").Value = -3822.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 323.95456
$ws.Range("I97").Value = 334.6842
$ws.Range("J97").Value = 256
$ws.Range("K97").Value = 334.6842
$ws.Range("L97").Value = 256
$ws.Range("M97").Value = 161.3158
$ws.Range("N97").Value = -1248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2184.9375
$ws.Range("I116").Value = 895.9
$ws.Range("J116").Value = 4333.3335
$ws.Range("K116").Value = 895.9
$ws.Range("L116").Value = 4333.3335
$ws.Range("M116").Value = 1398.1
$ws.Range("N116").Value = -8921.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4037.5
$ws.Range("I132").Value = 4037.5
$ws.Range("K132").Value = 12112.5
$ws.Range("M132").Value = -9582.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4395.7666
$ws.Range("I136").Value = 4034.8462
$ws.Range("K136").Value = 12104.5386
$ws.Range("M136").Value = -9554.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2184.9375
$ws.Range("I3").Value = 895.9
$ws.Range("J3").Value = 4333.3335
$ws.Range("K3").Value = 895.9
$ws.Range("L3").Value = 4333.3335
$ws.Range("M3").Value = -781.9
$ws.Range("N3").Value = -4561.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3581.3333
$ws.Range("I20").Value = 3454
$ws.Range("K20").Value = 3454
$ws.Range("M20").Value = -3207

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 869.1667
$ws.Range("I25").Value = 549.75
$ws.Range("K25").Value = 549.75
$ws.Range("M25").Value = -314.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2293.75
$ws.Range("I94").Value = 2150.05
$ws.Range("K94").Value = 2150.05
$ws.Range("M94").Value = -1699.05

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2519.76
$ws.Range("I99").Value = 2067.25
$ws.Range("K99").Value = 2067.25
$ws.Range("M99").Value = -569.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4052
$ws.Range("I134").Value = 4052
$ws.Range("K134").Value = 12156
$ws.Range("M134").Value = -9621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2237.5
$ws.Range("I12").Value = 700
$ws.Range("J12").Value = 2750
$ws.Range("K12").Value = 700
$ws.Range("L12").Value = 2750
$ws.Range("M12").Value = -530
$ws.Range("N12").Value = -3090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 68333.336
$ws.Range("J100").Value = 68333.336
$ws.Range("L100").Value = 68333.336
$ws.Range("N100").Value = -70497.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5710.8
$ws.Range("I125").Value = 5710.8
$ws.Range("K125").Value = 17132.4
$ws.Range("M125").Value = -12212.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20837802
$ws.Range("I70").Value = 37041380
$ws.Range("K70").Value = 37041380
$ws.Range("M70").Value = -37041110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 20837802
$ws.Range("I73").Value = 37041380
$ws.Range("K73").Value = 37041380
$ws.Range("M73").Value = -37040444

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3401.0715
$ws.Range("I22").Value = 2853.75
$ws.Range("J22").Value = 3811.5625
$ws.Range("K22").Value = 2853.75
$ws.Range("L22").Value = 3811.5625
$ws.Range("M22").Value = -2558.75
$ws.Range("N22").Value = -4401.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3401.0715
$ws.Range("I27").Value = 2853.75
$ws.Range("J27").Value = 3811.5625
$ws.Range("K27").Value = 2853.75
$ws.Range("L27").Value = 3811.5625
$ws.Range("M27").Value = -2746.75
$ws.Range("N27").Value = -4025.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4833.5
$ws.Range("I93").Value = 799
$ws.Range("J93").Value = 5640.4
$ws.Range("K93").Value = 799
$ws.Range("L93").Value = 5640.4
$ws.Range("M93").Value = 449
$ws.Range("N93").Value = -8136.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4309.8335
$ws.Range("I136").Value = 3270.5715
$ws.Range("K136").Value = 9811.7145
$ws.Range("M136").Value = -7261.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2237.375
$ws.Range("I81").Value = 2371.2856
$ws.Range("J81").Value = 1300
$ws.Range("K81").Value = 4742.5712
$ws.Range("L81").Value = 2600
$ws.Range("M81").Value = -3681.5712
$ws.Range("N81").Value = -4722

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2237.375
$ws.Range("I84").Value = 2371.2856
$ws.Range("J84").Value = 1300
$ws.Range("K84").Value = 23712.856
$ws.Range("L84").Value = 13000
$ws.Range("M84").Value = -18408.856
$ws.Range("N84").Value = -23608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4549.457
$ws.Range("I132").Value = 4647.839
$ws.Range("J132").Value = 3787
$ws.Range("K132").Value = 13943.517
$ws.Range("L132").Value = 11361
$ws.Range("M132").Value = -11413.517
$ws.Range("N132").Value = -16421

Write-Output "Applied 35 row updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
